# Update automatico via Actualizar 10-14-2020 17-33-01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily "Dolar observado" rows appended to the existing A1:B646 series,
# extending it through row 655 (dates 2020-10-06 .. 2020-10-14).
$newRows = @(
    @{ Row = 647; Date = [datetime]"2020-10-06"; Value = 794.34 },
    @{ Row = 648; Date = [datetime]"2020-10-07"; Value = 797.35 },
    @{ Row = 649; Date = [datetime]"2020-10-08"; Value = 795.05 },
    @{ Row = 650; Date = [datetime]"2020-10-09"; Value = 797.25 },
    @{ Row = 651; Date = [datetime]"2020-10-10"; Value = $null },
    @{ Row = 652; Date = [datetime]"2020-10-11"; Value = $null },
    @{ Row = 653; Date = [datetime]"2020-10-12"; Value = $null },
    @{ Row = 654; Date = [datetime]"2020-10-13"; Value = 796.05 },
    @{ Row = 655; Date = [datetime]"2020-10-14"; Value = 797.66 }
)

# Known good style donors already present in the sheet:
#  - row 646 col A/B: date style + numeric "observado" style
#  - row 645 col B:    "--" placeholder (no-market-day) style
$numericSrcRow = 646
$dashSrcRow = 645

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Seed formatting by copying from an existing donor row (same column
    # styles used throughout the existing date/value series), then
    # overwrite values only - this reuses existing style indices instead of
    # minting new ones.
    $dateSrc = $ws.Range("A" + $numericSrcRow)
    $dateDst = $ws.Range("A" + $rowNum)
    $dateSrc.Copy() | Out-Null
    $dateDst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    if ($null -eq $r.Value) {
        $valueSrc = $ws.Range("B" + $dashSrcRow)
    } else {
        $valueSrc = $ws.Range("B" + $numericSrcRow)
    }
    $valueDst = $ws.Range("B" + $rowNum)
    $valueSrc.Copy() | Out-Null
    $valueDst.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $dateCell = $ws.Cells.Item($rowNum, 1)
    $dateCell.Value = $r.Date

    $valueCell = $ws.Cells.Item($rowNum, 2)
    if ($null -eq $r.Value) {
        $valueCell.Value = "--"
    } else {
        $valueCell.Value = $r.Value
    }
}

# Extend the named range and refresh the view to match the new extent.
$wb.Names.Item("DOLAR_OBS_ADO").RefersTo = "=DOLAR_OBS_ADO!`$A`$1:`$B`$655"

$ws.Application.ActiveWindow.ScrollRow = 649
$ws.Range("B658").Select() | Out-Null
